# Update column G ("K") values on the active worksheet to reflect the
# regenerated save_data (switch from Strike# to K, recalculated s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 4
    4  = 7
    5  = 7
    6  = 4
    7  = 7
    8  = 8
    9  = 8
    10 = 4
    11 = 4
    12 = 2
    13 = 7
    14 = 2
    15 = 2
    16 = 5
    17 = 3
    18 = 9
    19 = 4
    20 = 4
    21 = 2
    22 = 3
    23 = 7
    24 = 2
    25 = 8
    26 = 5
    27 = 5
    28 = 1
    29 = 3
    30 = 10
    31 = 6
    32 = 1
    33 = 8
    34 = 2
    35 = 2
    36 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
